# Insert a new weekly price-report row for "Femacal de La Calera - Bruselas
# (repollito)" right after the existing row 33, shifting rows 34:59 down to
# 35:60 (dimension grows from A1:R59 to A1:R60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(34).Insert()

$ws.Range("A34").Value = 3
$ws.Range("B34").Value = "Femacal de La Calera"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 44762
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = 100112035
$ws.Range("G34").Value = "Bruselas (repollito)"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 80
$ws.Range("K34").Value = 14000
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = 14500
$ws.Range("N34").Value = "$/malla 15 kilos"
$ws.Range("O34").Value = "Provincia de Quillota"
$ws.Range("P34").Value = 967
$ws.Range("Q34").Value = 15
$ws.Range("R34").Value = "Hortaliza"
